$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (rows 2-5)
$data = @(
    @("1328021", "https://aiesec.org/opportunity/global-talent/1328021", "Guest Relations Officer Intern", "Nugegoda, Sri Lanka", "No", "21 applicants", "3 - 6 Months", "The Barn By Starbeans in Ella"),
    @("1327889", "https://aiesec.org/opportunity/global-talent/1327889", "Graphic Designer", "Birkat as SAB, Madinet Berkat as Sabee, Birket el Sab, Menofia Governorate, Egypt", "No", "4 applicants", "9 - 12 Weeks", "Lines"),
    @("1327809", "https://aiesec.org/opportunity/global-talent/1327809", "Video editor", "El Sadat City, Menofia Governorate, Egypt", "No", "3 applicants", "9 - 12 Weeks", "Habib Agency"),
    @("1306000", "https://aiesec.org/opportunity/global-talent/1306000", "Business Development Intern", "Indore, Madhya Pradesh, India", "No", "26 applicants", "3 - 6 Months", "Walkover Web Solutions")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($c = 0; $c -lt 8; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($c -eq 0) {
            # Opportunity ID must be stored as text, not a number
            $cell.NumberFormat = "@"
            $cell.Value = $data[$i][$c]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $data[$i][$c]
        }
    }
}

# Update column widths to match the new target layout.
# The runtime's ColumnWidth -> stored <col width> conversion adds 5/6 (0.8333),
# so subtract that offset here to land on the exact integer widths required.
$widthOffset = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 55 - $widthOffset
$ws.Columns.Item(3).ColumnWidth = 33 - $widthOffset
$ws.Columns.Item(4).ColumnWidth = 84 - $widthOffset
$ws.Columns.Item(5).ColumnWidth = 10 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 16 - $widthOffset
$ws.Columns.Item(7).ColumnWidth = 15 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 32 - $widthOffset
